# US-2.2_Accountant_Account receivable functionality_TCs.xlsx
#
# The ID fields (receivableID / entryID / customerID) that show up on the
# "Validate the account receivable information page" (TC02) and on the
# "Error validation when adding new accounts receivable enteries" (TC07)
# test cases were wrongly documented as fields that "should be displayed".
# These are internal/system generated identifiers, so the expected result
# is corrected to state that they should NOT be displayed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "It should not be displayed and should be as per parameter."

# TC02 - "Validate the account receivable information page."
#   Step 12 -> Validate receivableID field.
#   Step 13 -> Validate entryID field.
#   Step 14 -> Validate customerID field.
$ws.Range("F19").Value = $newText
$ws.Range("F20").Value = $newText
$ws.Range("F21").Value = $newText

# TC(duplicate block) - same three "Validate …ID field." steps repeated
# further down the sheet for the "add new accounts receivable enteries"
# test case.
$ws.Range("F46").Value = $newText
$ws.Range("F47").Value = $newText
$ws.Range("F48").Value = $newText

# TC07 - "Error validation when adding new accounts receivable enteries"
#   Step 11 -> Validate receivableID field.
#   Step 12 -> Validate entryID field.
#   Step 13 -> Validate customerID field.
$ws.Range("F68").Value = $newText
$ws.Range("F69").Value = $newText
$ws.Range("F70").Value = $newText

# The new wording is shorter than the old one, so these three rows no
# longer need the extra wrapped-text height and shrink back down to the
# sheet's default row height (row 71 keeps its old text/height).
$ws.Rows.Item(68).AutoFit()
$ws.Rows.Item(69).AutoFit()
$ws.Rows.Item(70).AutoFit()

# Re-position the view the way the author left it when saving: scrolled
# down so row 58 is the top visible row, with E77 as the active/selected
# cell.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 58
$ws.Range("E77").Select()
